# The workbook is an events/agenda schedule on Sheet1 (A1:F9).
# This script updates session/activity descriptions per the new agenda text,
# turns the markdown-link Session name in the "Activity 1" row back into
# plain text, re-applies wrap-text + row-height formatting to the rows whose
# Description text grew, and updates the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Description (column E) text for several sessions ---
$ws.Range("E2").Value = "Greetings from the FAMPS and FSN Chairs; Highlights from Day 1"
$ws.Range("E3").Value = "Challenges and bottlenecks of working with administrative data"

# --- Row 6 ("Activity 1"): Session cell was a markdown link, now plain text ---
$ws.Range("D6").Value = "Activity 1"
$ws.Range("E6").Value = "Navigating Licenses Building a Research Plan to Access RDC Data"

$ws.Range("E7").Value = "Deterministic Data Linkages: Matching and Fuzzy Matching"
$ws.Range("E8").Value = "Econometric and Statistical Considerations When Using Linked Data"
$ws.Range("E9").Value = "Closing from the FAMPS and FSN Chairs"

# --- Wrap text for the Description cells that now hold longer text ---
$ws.Range("E6").WrapText = $true
$ws.Range("E7").WrapText = $true
$ws.Range("E8").WrapText = $true

# --- Row height adjustments to fit the new wrapped text ---
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 30

# --- Update the saved selection/active cell ---
$ws.Range("F12").Select()
